$d = $word.ActiveDocument

# 1. Strike through the "Use padding instead of nbsp's for title boxes" bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Use padding instead of*nbsp*for title boxes*") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# 2. Merge the "Update screenshots of " + "devices" runs into a single run.
$range = $d.Content
$range.Find.Execute("Update screenshots of devices", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "Update screenshots of devices", 2)
